$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 30.914286
$ws.Range("I4").Value = 30.914286
$ws.Range("K4").Value = 30.914286
$ws.Range("M4").Value = 83.085714

$ws.Range("H6").Value = 56.3125
$ws.Range("I6").Value = 56.3125
$ws.Range("K6").Value = 168.9375
$ws.Range("M6").Value = -56.9375

$ws.Range("H15").Value = 890.87756
$ws.Range("I15").Value = 890.87756
$ws.Range("K15").Value = 2672.63268
$ws.Range("M15").Value = -2503.63268

$ws.Range("H43").Value = 4228.93
$ws.Range("I43").Value = 4235.896
$ws.Range("K43").Value = 4235.896
$ws.Range("M43").Value = -4166.896

$ws.Range("H46").Value = 2739.5
$ws.Range("J46").Value = 3000
$ws.Range("L46").Value = 9000
$ws.Range("N46").Value = -9238

$ws.Range("H60").Value = 2739.5
$ws.Range("J60").Value = 3000
$ws.Range("L60").Value = 9000
$ws.Range("N60").Value = -9968

$ws.Range("H80").Value = 871.675
$ws.Range("I80").Value = 579.05
$ws.Range("J80").Value = 1164.3
$ws.Range("K80").Value = 1737.15
$ws.Range("L80").Value = 3492.9
$ws.Range("M80").Value = -739.1499999999999
$ws.Range("N80").Value = -5488.9

$ws.Range("H83").Value = 871.675
$ws.Range("I83").Value = 579.05
$ws.Range("J83").Value = 1164.3
$ws.Range("K83").Value = 5211.45
$ws.Range("L83").Value = 10478.7
$ws.Range("M83").Value = -219.4499999999998
$ws.Range("N83").Value = -20462.7

$ws.Range("H92").Value = 4066.4092
$ws.Range("I92").Value = 4529.4375
$ws.Range("J92").Value = 2831.6667
$ws.Range("K92").Value = 4529.4375
$ws.Range("L92").Value = 2831.6667
$ws.Range("M92").Value = -3281.4375
$ws.Range("N92").Value = -5327.6667

$ws.Range("H129").Value = 17901.316
$ws.Range("I129").Value = 47108.934
$ws.Range("J129").Value = 1050.7693
$ws.Range("K129").Value = 141326.802
$ws.Range("L129").Value = 3152.3079
$ws.Range("M129").Value = -136326.802
$ws.Range("N129").Value = -13152.3079

$ws.Range("H135").Value = 3484.2173
$ws.Range("I135").Value = 1459.35
$ws.Range("J135").Value = 16983.334
$ws.Range("K135").Value = 13134.15
$ws.Range("L135").Value = 152850.006
$ws.Range("M135").Value = -10599.15
$ws.Range("N135").Value = -157920.006

$ws.Range("H137").Value = 11627.542
$ws.Range("I137").Value = 4035.0625
$ws.Range("J137").Value = 26812.5
$ws.Range("K137").Value = 12105.1875
$ws.Range("L137").Value = 80437.5
$ws.Range("M137").Value = -9555.1875
$ws.Range("N137").Value = -85537.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5764.295
$ws.Range("I32").Value = 3095.758
$ws.Range("K32").Value = 3095.758
$ws.Range("M32").Value = -2808.758

$ws.Range("H51").Value = 39790.332
$ws.Range("J51").Value = 39790.332
$ws.Range("L51").Value = 39790.332
$ws.Range("N51").Value = -41302.332

$ws.Range("H61").Value = 14213.107
$ws.Range("I61").Value = 2753.389
$ws.Range("J61").Value = 34840.6
$ws.Range("K61").Value = 2753.389
$ws.Range("L61").Value = 34840.6
$ws.Range("M61").Value = -2541.389
$ws.Range("N61").Value = -35264.6

$ws.Range("H136").Value = 14213.107
$ws.Range("I136").Value = 2753.389
$ws.Range("J136").Value = 34840.6
$ws.Range("K136").Value = 8260.167000000001
$ws.Range("L136").Value = 104521.8
$ws.Range("M136").Value = -5710.167000000001
$ws.Range("N136").Value = -109621.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 14553.5
$ws.Range("J99").Value = 21645.8
$ws.Range("L99").Value = 21645.8
$ws.Range("N99").Value = -24641.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 10761.917
$ws.Range("I16").Value = 5293.3076
$ws.Range("K16").Value = 5293.3076
$ws.Range("M16").Value = -5006.3076

$ws.Range("H29").Value = 3350
$ws.Range("J29").Value = 6000
$ws.Range("L29").Value = 6000
$ws.Range("N29").Value = -6586

$ws.Range("H105").Value = 34536.668
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()

$ws.Range("H113").Value = 10761.917
$ws.Range("I113").Value = 5293.3076
$ws.Range("K113").Value = 5293.3076
$ws.Range("M113").Value = -3123.3076

$ws.Range("H132").Value = 4890.4546
$ws.Range("I132").Value = 1875.52
$ws.Range("K132").Value = 5626.559999999999
$ws.Range("M132").Value = -3096.559999999999

$ws.Range("H134").Value = 28577218
$ws.Range("I134").Value = 1419
$ws.Range("J134").Value = 47627750
$ws.Range("K134").Value = 4257
$ws.Range("L134").Value = 142883250
$ws.Range("M134").Value = -1722
$ws.Range("N134").Value = -142888320

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("I5").Value = 1048.1852
$ws.Range("K5").Value = 3144.5556
$ws.Range("M5").Value = -3032.5556

$ws.Range("H51").Value = 3000
$ws.Range("J51").Value = 3000
$ws.Range("L51").Value = 9000
$ws.Range("N51").Value = -9920

$ws.Range("H94").Value = 15451
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 15451
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 46353
$ws.Range("M94").ClearContents()
$ws.Range("N94").Value = -47705

$ws.Range("H107").Value = 921378.2
$ws.Range("I107").Value = 495.8
$ws.Range("J107").Value = 1648390.6
$ws.Range("K107").Value = 1487.4
$ws.Range("L107").Value = 4945171.800000001
$ws.Range("M107").Value = 432.5999999999999
$ws.Range("N107").Value = -4949011.800000001

$ws.Range("H113").Value = 861.7406999999999
$ws.Range("I113").Value = 351.36365
$ws.Range("K113").Value = 1054.09095
$ws.Range("M113").Value = 1115.90905

$ws.Range("H131").Value = 1495.6162
$ws.Range("J131").Value = 1495.6162
$ws.Range("L131").Value = 4486.848599999999
$ws.Range("N131").Value = -14566.8486

$ws.Range("I135").Value = 1048.1852
$ws.Range("K135").Value = 9433.666799999999
$ws.Range("M135").Value = -6898.666799999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 14178.64
$ws.Range("I80").Value = 16894.273
$ws.Range("J80").Value = 12044.929
$ws.Range("K80").Value = 16894.273
$ws.Range("L80").Value = 12044.929
$ws.Range("M80").Value = -15896.273
$ws.Range("N80").Value = -14040.929

$ws.Range("H83").Value = 14178.64
$ws.Range("I83").Value = 16894.273
$ws.Range("J83").Value = 12044.929
$ws.Range("K83").Value = 84471.36500000001
$ws.Range("L83").Value = 60224.645
$ws.Range("M83").Value = -79479.36500000001
$ws.Range("N83").Value = -70208.645

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 20002730
$ws.Range("J22").Value = 66669580
$ws.Range("L22").Value = 66669580
$ws.Range("N22").Value = -66670170

$ws.Range("H27").Value = 20002730
$ws.Range("J27").Value = 66669580
$ws.Range("L27").Value = 66669580
$ws.Range("N27").Value = -66669794

$ws.Range("H105").Value = 30000
$ws.Range("J105").Value = 30000
$ws.Range("L105").Value = 30000
$ws.Range("N105").Value = -36988

$ws.Range("H122").Value = 27783160
$ws.Range("I122").Value = 32262734
$ws.Range("K122").Value = 96788202
$ws.Range("M122").Value = -96785752

$ws.Range("H132").Value = 1679480.6
$ws.Range("I132").Value = 5201.25
$ws.Range("J132").Value = 10050877
$ws.Range("K132").Value = 15603.75
$ws.Range("L132").Value = 30152631
$ws.Range("M132").Value = -13073.75
$ws.Range("N132").Value = -30157691

$ws.Range("H136").Value = 16348.12
$ws.Range("I136").Value = 16635.357
$ws.Range("J136").Value = 15982.546
$ws.Range("K136").Value = 49906.071
$ws.Range("L136").Value = 47947.638
$ws.Range("M136").Value = -47356.071
$ws.Range("N136").Value = -53047.638

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2109.2
$ws.Range("J96").Value = 2190
$ws.Range("L96").Value = 2190
$ws.Range("N96").Value = -4936

$ws.Range("H122").Value = 929409.75
$ws.Range("I122").Value = 1308402
$ws.Range("J122").Value = 9000
$ws.Range("K122").Value = 3925206
$ws.Range("L122").Value = 27000
$ws.Range("M122").Value = -3922756
$ws.Range("N122").Value = -31900

$ws.Range("H132").Value = 6063.2324
$ws.Range("I132").Value = 2809.5151
$ws.Range("K132").Value = 8428.5453
$ws.Range("M132").Value = -5898.5453

Write-Host "Updated Leve profit figures across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets"
